# THE GHOUL.docx -- apply story-text revisions described in the commit.
$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $findText"
    }
}

# "...kingdom name Oak. The Oak has..." -> "...kingdom named Oak. Oak has..."
Replace-Text "kingdom name Oak. The Oak has" "kingdom named Oak. Oak has"

# "...protect them from monsters..." -> "...protected from monsters..."
Replace-Text "protect them" "protected"

# "...gave life to crops and nature covered. The..." -> "...gave life to crops and nature. The..."
Replace-Text "nature covered. The" "nature. The"

# Wisdom Stone sentence rewrite
Replace-Text "enhance the intuition of the awareness, and pass the extensive knowledge to the holder" "enhanced intuition, awareness, and passed extensive knowledge of ancestors to the holder"

# "King Luis first" -> "King Luis I (the first)"
Replace-Text "King Luis first" "King Luis I (the first)"

# Oath sentence rewrite
Replace-Text ". After that, they all swear to the King that they and their family will protect the stone forever. " ". All holders vowed under oath to have their lineage protect the stones forever…"

# "Two hundred years later, the third king Alex, ascends the throne. One day, " -> "Two hundred years later King Alex ascends the throne. One day, "
Replace-Text "Two hundred years later, the third king Alex, ascends the throne. One day, " "Two hundred years later King Alex ascends the throne. One day, "

# "The General was missing on his way" -> "The General went missing on his way"
Replace-Text "The General was missing on his way" "The General went missing on his way"

# Famine/population sentence rewrite
Replace-Text "More and more hungry people become homeless and get attacked by the Ghoul." "The population rapidly declined due to famine and attacks from the Ghoul."

# "You are the Guard knights who own the wisdom stone. " -> "You are a member of the Guard knights who has been entrusted with the wisdom stone. "
Replace-Text "You are the Guard knights who own the wisdom stone. " "You are a member of the Guard knights who has been entrusted with the wisdom stone. "

# Remove the now-duplicate blank "NoSpacing" paragraph that separated the
# famine paragraph from the "You are..." paragraph (two blank paragraphs
# collapse into one).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "`r" -and $para.Range.Text.Trim() -eq "") {
        $next = $para.Next()
        if ($next -ne $null -and $next.Range.Text.StartsWith("You are ")) {
            $para.Range.Delete()
            break
        }
    }
}
